$d = $word.ActiveDocument

function Set-ParagraphXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $countBefore = $d.Paragraphs.Count
    $pkgHead = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">'
    $pkgTail = '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $pkg = $pkgHead + $innerXml + $pkgTail
    $p.Range.InsertXML($pkg)
    $countAfter = $d.Paragraphs.Count
    if ($countAfter -gt $countBefore) {
        # InsertXML against the very last paragraph's range (which includes
        # the trailing paragraph mark) leaves a spurious empty paragraph
        # behind after the newly-inserted one; collapse it back out.
        $extraCount = $countAfter - $countBefore
        $k = 0
        while ($k -lt $extraCount) {
            $newPara = $d.Paragraphs.Item($paraIndex)
            $markStart = $newPara.Range.End - 1
            $markEnd = $newPara.Range.End
            $markRange = $d.Range($markStart, $markEnd)
            $markRange.Delete()
            $k = $k + 1
        }
    }
}

# 1. "CS1501 Nodejs" -> "CS1501 " + "Nodejs" (two runs, spell-check proof markers)
$xml2 = '<w:body><w:p w14:paraId="15CB6328" w14:textId="77777777" w:rsidR="00EC6C6E" w:rsidRPr="006D1818" w:rsidRDefault="006D1818" w:rsidP="006D1818">'
$xml2 = $xml2 + '<w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Source Sans Pro" w:hAnsi="Source Sans Pro"/><w:sz w:val="96"/><w:szCs w:val="96"/></w:rPr></w:pPr>'
$xml2 = $xml2 + '<w:r><w:rPr><w:rFonts w:ascii="Source Sans Pro" w:hAnsi="Source Sans Pro"/><w:sz w:val="96"/><w:szCs w:val="96"/></w:rPr><w:t xml:space="preserve">CS1501 </w:t></w:r>'
$xml2 = $xml2 + '<w:proofErr w:type="spellStart"/>'
$xml2 = $xml2 + '<w:r><w:rPr><w:rFonts w:ascii="Source Sans Pro" w:hAnsi="Source Sans Pro"/><w:sz w:val="96"/><w:szCs w:val="96"/></w:rPr><w:t>Nodejs</w:t></w:r>'
$xml2 = $xml2 + '<w:proofErr w:type="spellEnd"/>'
$xml2 = $xml2 + '</w:p></w:body>'
$idx2 = 2
Set-ParagraphXml $idx2 $xml2

# 2. "Sublime Text Tips and Tricks" -> "HTML and CSS"
$xml3 = '<w:body><w:p w14:paraId="36A7AEA1" w14:textId="77777777" w:rsidR="006D1818" w:rsidRDefault="002E5C61" w:rsidP="006D1818">'
$xml3 = $xml3 + '<w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Source Sans Pro" w:hAnsi="Source Sans Pro"/><w:sz w:val="96"/><w:szCs w:val="96"/></w:rPr></w:pPr>'
$xml3 = $xml3 + '<w:r><w:rPr><w:rFonts w:ascii="Source Sans Pro" w:hAnsi="Source Sans Pro"/><w:sz w:val="96"/><w:szCs w:val="96"/></w:rPr><w:t>HTML and CSS</w:t></w:r>'
$xml3 = $xml3 + '</w:p></w:body>'
$idx3 = 3
Set-ParagraphXml $idx3 $xml3

# 3. "Feb 18" -> "Feb 24" (keep ", 2014" as its own separate run)
$xml5 = '<w:body><w:p w14:paraId="785DBF3E" w14:textId="77777777" w:rsidR="006D1818" w:rsidRDefault="002E5C61" w:rsidP="006D1818">'
$xml5 = $xml5 + '<w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Source Sans Pro" w:hAnsi="Source Sans Pro"/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr></w:pPr>'
$xml5 = $xml5 + '<w:r><w:rPr><w:rFonts w:ascii="Source Sans Pro" w:hAnsi="Source Sans Pro"/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>Feb 24</w:t></w:r>'
$xml5 = $xml5 + '<w:r w:rsidR="006D1818"><w:rPr><w:rFonts w:ascii="Source Sans Pro" w:hAnsi="Source Sans Pro"/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>, 2014</w:t></w:r>'
$xml5 = $xml5 + '</w:p></w:body>'
$idx5 = 5
Set-ParagraphXml $idx5 $xml5

# 4. "Get to know sublime text" -> "Introduction to HTML, CSS, and Bootstrap"
#    (split across the existing _GoBack bookmark into two runs: "...Bootstra" + "p")
$xml7 = '<w:body><w:p w14:paraId="76943014" w14:textId="77777777" w:rsidR="006D1818" w:rsidRPr="006D1818" w:rsidRDefault="002E5C61" w:rsidP="006D1818">'
$xml7 = $xml7 + '<w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Source Sans Pro" w:hAnsi="Source Sans Pro"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>'
$xml7 = $xml7 + '<w:r><w:rPr><w:rFonts w:ascii="Source Sans Pro" w:hAnsi="Source Sans Pro"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Introduction to HTML, CSS, and Bootstra</w:t></w:r>'
$xml7 = $xml7 + '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$xml7 = $xml7 + '<w:r><w:rPr><w:rFonts w:ascii="Source Sans Pro" w:hAnsi="Source Sans Pro"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>p</w:t></w:r>'
$xml7 = $xml7 + '</w:p></w:body>'
$idx7 = 7
Set-ParagraphXml $idx7 $xml7
